# Fix bug: exceeded request in google drive
# Update price list date and the "SOGA de Monofilamento" revestida prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Date stamp in A1 moves forward one day (45310 -> 45311)
$ws.Range("A1").Value = 45311

# Updated prices x metro for the revestida soga rolls
$ws.Range("D14").Value = 43.2
$ws.Range("D15").Value = 61.8
$ws.Range("D16").Value = 88.40000000000001
$ws.Range("D17").Value = 154
